$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 20:19:52"
$ws1.Cells.Item(3,1).Value = "Total filas: 615"

$ws1.Cells.Item(600,2).Value = "20:19:41"
$ws1.Cells.Item(600,3).Value = "20:22"
$ws1.Cells.Item(600,4).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(600,5).Value = 3
$ws1.Cells.Item(600,6).Value = "LP1912"
$ws1.Cells.Item(600,7).Value = "30/12/2025"

$ws1.Cells.Item(601,2).Value = "20:19:41"
$ws1.Cells.Item(601,3).Value = "20:22"
$ws1.Cells.Item(601,4).Value = "16_SANTA ANA"
$ws1.Cells.Item(601,5).Value = 3
$ws1.Cells.Item(601,6).Value = "LP1912"
$ws1.Cells.Item(601,7).Value = "30/12/2025"

$ws1.Cells.Item(602,2).Value = "20:19:41"
$ws1.Cells.Item(602,3).Value = "20:23"
$ws1.Cells.Item(602,4).Value = "215A_EL PATO"
$ws1.Cells.Item(602,5).Value = 4
$ws1.Cells.Item(602,6).Value = "LP1912"
$ws1.Cells.Item(602,7).Value = "30/12/2025"

$ws1.Cells.Item(603,2).Value = "20:19:41"
$ws1.Cells.Item(603,3).Value = "20:34"
$ws1.Cells.Item(603,4).Value = "16_SANTA ANA"
$ws1.Cells.Item(603,5).Value = 15
$ws1.Cells.Item(603,6).Value = "LP1912"
$ws1.Cells.Item(603,7).Value = "30/12/2025"

$ws1.Cells.Item(604,2).Value = "20:19:41"
$ws1.Cells.Item(604,3).Value = "20:45"
$ws1.Cells.Item(604,4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(604,5).Value = 26
$ws1.Cells.Item(604,6).Value = "LP1912"
$ws1.Cells.Item(604,7).Value = "30/12/2025"

$ws1.Cells.Item(605,2).Value = "20:19:41"
$ws1.Cells.Item(605,3).Value = "20:46"
$ws1.Cells.Item(605,4).Value = "16_SANTA ANA"
$ws1.Cells.Item(605,5).Value = 27
$ws1.Cells.Item(605,6).Value = "LP1912"
$ws1.Cells.Item(605,7).Value = "30/12/2025"

$ws1.Cells.Item(606,2).Value = "20:19:41"
$ws1.Cells.Item(606,3).Value = "20:52"
$ws1.Cells.Item(606,4).Value = "15_ABASTO"
$ws1.Cells.Item(606,5).Value = 33
$ws1.Cells.Item(606,6).Value = "LP1912"
$ws1.Cells.Item(606,7).Value = "30/12/2025"

$ws1.Cells.Item(607,2).Value = "20:19:41"
$ws1.Cells.Item(607,3).Value = "20:57"
$ws1.Cells.Item(607,4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(607,5).Value = 38
$ws1.Cells.Item(607,6).Value = "LP1912"
$ws1.Cells.Item(607,7).Value = "30/12/2025"

$ws1.Cells.Item(608,2).Value = "20:19:41"
$ws1.Cells.Item(608,3).Value = "21:04"
$ws1.Cells.Item(608,4).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(608,5).Value = 45
$ws1.Cells.Item(608,6).Value = "LP1912"
$ws1.Cells.Item(608,7).Value = "30/12/2025"

$ws1.Cells.Item(609,2).Value = "20:19:41"
$ws1.Cells.Item(609,3).Value = "21:07"
$ws1.Cells.Item(609,4).Value = "215B_EL PATO"
$ws1.Cells.Item(609,5).Value = 48
$ws1.Cells.Item(609,6).Value = "LP1912"
$ws1.Cells.Item(609,7).Value = "30/12/2025"

$ws1.Cells.Item(610,2).Value = "20:19:41"
$ws1.Cells.Item(610,3).Value = "21:20"
$ws1.Cells.Item(610,4).Value = "26_HERNANDEZ"
$ws1.Cells.Item(610,5).Value = 61
$ws1.Cells.Item(610,6).Value = "LP1912"
$ws1.Cells.Item(610,7).Value = "30/12/2025"

$ws1.Cells.Item(611,2).Value = "20:19:41"
$ws1.Cells.Item(611,3).Value = "21:22"
$ws1.Cells.Item(611,4).Value = "15_ABASTO"
$ws1.Cells.Item(611,5).Value = 63
$ws1.Cells.Item(611,6).Value = "LP1912"
$ws1.Cells.Item(611,7).Value = "30/12/2025"

$ws1.Cells.Item(612,2).Value = "20:19:41"
$ws1.Cells.Item(612,3).Value = "21:32"
$ws1.Cells.Item(612,4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(612,5).Value = 73
$ws1.Cells.Item(612,6).Value = "LP1912"
$ws1.Cells.Item(612,7).Value = "30/12/2025"

$ws1.Cells.Item(613,2).Value = "20:19:41"
$ws1.Cells.Item(613,3).Value = "21:37"
$ws1.Cells.Item(613,4).Value = "17_ROMERO"
$ws1.Cells.Item(613,5).Value = 78
$ws1.Cells.Item(613,6).Value = "LP1912"
$ws1.Cells.Item(613,7).Value = "30/12/2025"

$ws1.Cells.Item(614,2).Value = "20:19:41"
$ws1.Cells.Item(614,3).Value = "21:42"
$ws1.Cells.Item(614,4).Value = "14_ABASTO"
$ws1.Cells.Item(614,5).Value = 83
$ws1.Cells.Item(614,6).Value = "LP1912"
$ws1.Cells.Item(614,7).Value = "30/12/2025"

$ws1.Cells.Item(615,2).Value = "20:19:41"
$ws1.Cells.Item(615,3).Value = "21:47"
$ws1.Cells.Item(615,4).Value = "215A_EL PATO"
$ws1.Cells.Item(615,5).Value = 88
$ws1.Cells.Item(615,6).Value = "LP1912"
$ws1.Cells.Item(615,7).Value = "30/12/2025"

$ws1.Cells.Item(616,2).Value = "20:19:41"
$ws1.Cells.Item(616,3).Value = "21:51"
$ws1.Cells.Item(616,4).Value = "10_OLMOS"
$ws1.Cells.Item(616,5).Value = 92
$ws1.Cells.Item(616,6).Value = "LP1912"
$ws1.Cells.Item(616,7).Value = "30/12/2025"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 20:19:52"
$ws2.Cells.Item(3,1).Value = "Total filas: 44"

$ws2.Cells.Item(43,2).Value = "30/12/2025"
$ws2.Cells.Item(43,3).Value = "20:19:41"
$ws2.Cells.Item(43,4).Value = "20:23"
$ws2.Cells.Item(43,5).Value = "215A_EL PATO"
$ws2.Cells.Item(43,6).Value = 4
$ws2.Cells.Item(43,7).Value = "LP1912"

$ws2.Cells.Item(44,2).Value = "30/12/2025"
$ws2.Cells.Item(44,3).Value = "20:19:41"
$ws2.Cells.Item(44,4).Value = "21:07"
$ws2.Cells.Item(44,5).Value = "215B_EL PATO"
$ws2.Cells.Item(44,6).Value = 48
$ws2.Cells.Item(44,7).Value = "LP1912"

$ws2.Cells.Item(45,2).Value = "30/12/2025"
$ws2.Cells.Item(45,3).Value = "20:19:41"
$ws2.Cells.Item(45,4).Value = "21:47"
$ws2.Cells.Item(45,5).Value = "215A_EL PATO"
$ws2.Cells.Item(45,6).Value = 88
$ws2.Cells.Item(45,7).Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 20:19:52"
$ws3.Cells.Item(3,1).Value = "Total filas: 74"

$ws3.Cells.Item(75,2).Value = "30/12/2025"
$ws3.Cells.Item(75,3).Value = "20:19:46"
$ws3.Cells.Item(75,4).Value = "21:29"
$ws3.Cells.Item(75,5).Value = "215C_LA PLATA"
$ws3.Cells.Item(75,6).Value = 70
$ws3.Cells.Item(75,7).Value = "L6203"
